# Add the new accelerometer rows (22-31) captured on May 9th.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(20.3416576385498, -37.60380554199219, 9.662521362304688),
    @(-8.722187995910645, 8.096003532409668, -8.941717147827148),
    @(-72.57207489013672, -62.40628051757812, 63.55230712890625),
    @(-86.77053833007812, -16.04696655273438, -37.80954742431641),
    @(-12.0439281463623, -10.04275798797607, -17.28581237792969),
    @(-29.28538513183594, -21.46180725097656, 6.313263416290283),
    @(3.212208271026612, 6.285176277160645, 0.2726368904113769),
    @(-6.000577926635742, 12.46978282928467, 1.907239437103272),
    @(19.4264965057373, -45.61504745483398, 8.111544609069824),
    @(-43.24179840087891, 13.53559017181396, -9.444896697998049)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
